# Add three new artwork entries to the "DB" sheet / Tabla1 table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")
$ws.Activate()

$tbl = $ws.ListObjects.Item("Tabla1")

# Grow the table by three rows (Tabla1 currently spans A1:I146 -> A1:I149)
$newRow1 = $tbl.ListRows.Add()
$newRow2 = $tbl.ListRows.Add()
$newRow3 = $tbl.ListRows.Add()

$r1 = $newRow1.Range.Row
$r2 = $newRow2.Range.Row
$r3 = $newRow3.Range.Row

# Copy the formatting of the previous last data row onto the new rows so the
# date/centered styles match the rest of the table.
$ws.Range("A146:H146").Copy()
$ws.Range("A147:H149").PasteSpecial(-4122)

# Fill in the text columns (Nombre / Inspiracion) in this specific order so
# that the workbook's shared-string table grows with the expected entries.
$ws.Range("B147").Value = "Barranquilla"
$ws.Range("B149").Value = "Matterhorn"
$ws.Range("B148").Value = "Volcano"
$ws.Range("C148").Value = "Jocelyn Carmes"
$ws.Range("C147").Value = "Frederic Church"
$ws.Range("C149").Value = "Albert Bierstadt"

# Dates (Fecha)
$ws.Range("A147").Value = 45120
$ws.Range("A148").Value = 45139
$ws.Range("A149").Value = 45148

# Tecnica
$ws.Range("F147").Value = "Acuarela"
$ws.Range("F148").Value = "Oleo"
$ws.Range("F149").Value = "Oleo"

# Temática
$ws.Range("H147").Value = "Escena"
$ws.Range("H148").Value = "Escena"
$ws.Range("H149").Value = "Paisaje"

# These rows have no Coleccion (column G) entry, unlike row 146.
$ws.Range("G147").ClearContents()
$ws.Range("G148").ClearContents()
$ws.Range("G149").ClearContents()

# Leave the selection where data entry ended, matching the saved view state.
$ws.Range("H150").Select()
